$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper now also captures height/weight. Those two new columns are
# inserted at E/F, so the pre-existing "fantasy points" column (formerly E)
# moves to G. Shift that data over first, before writing the new columns.
for ($r = 2; $r -le 16; $r++) {
    $fantasyPoints = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 7).Value = $fantasyPoints
}

# Give the two new header cells (F1 currently blank, G1 will hold the
# relocated "fantasy points" header) the same look as the other header
# cells (bold, bordered, centered) by copying the format from D1.
$ws.Range("D1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# Fill in the new height/weight data for every data row
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.166666666666667
    $ws.Cells.Item($r, 6).Value = 245
}
